# Revisi kombinasi loc: refresh B:G metrics for AK..TD rows, re-label the
# FY/TO boundary rows (FY_4 absorbs old TO_1's slot, TO/TD rows shift down
# by one), and append the new TD_3/TD_4 rows (sheet grows to A1:G45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('AK_1.png', 0.354, 0.001, 0.496, 0.857, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_2.png', 0.6929999999999999, 0.002, 0.487, 1, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_3.png', 0.675, 0.002, 0.649, 1, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_4.png', 0.328, 0.001, 0.47, 1, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_5.png', 0.328, 0.001, 0.551, 1, 'Akhlak Kamiswara', 'Benar')
    ,@('MIB_1.png', 1.346, 0.004, 0.134, 0.714, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_2.png', 0.999, 0.003, 0.344, 1, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_3.png', 1.567, 0.005, 0.706, 1, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_4.png', 0.985, 0.003, 0.251, 1, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_5.png', 1.331, 0.004, 0.269, 1, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('AAH_1.png', 0.703, 0.002, 0.287, 1, 'Andrea Ayunove Hutami', 'Benar')
    ,@('AAH_2.png', 0.944, 0.003, 0.721, 1, 'Andrea Ayunove Hutami', 'Benar')
    ,@('AAH_3.png', 0.717, 0.002, 0.311, 1, 'Andrea Ayunove Hutami', 'Benar')
    ,@('TI_1.png', 0.804, 0.003, 0.502, 1, 'Toni Ismail', 'Benar')
    ,@('TI_2.png', 0.776, 0.003, 0.498, 1, 'Toni Ismail', 'Benar')
    ,@('TI_3.png', 0.551, 0.002, 0.795, 1, 'Toni Ismail', 'Benar')
    ,@('TI_4.png', 0.529, 0.002, 0.511, 1, 'Toni Ismail', 'Benar')
    ,@('TI_5.png', 0.875, 0.003, 0.537, 1, 'Toni Ismail', 'Benar')
    ,@('RAS_1.png', 0.455, 0.002, 0.498, 1, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RAS_2.png', 0.85, 0.003, 0.387, 1, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RAS_3.png', 0.378, 0.001, 0.383, 0.857, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RAS_4.png', 1.087, 0.004, 0.181, 0.714, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RAS_5.png', 1.01, 0.003, 0.502, 1, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RR_1.png', 1.302, 0.004, 0.5639999999999999, 1, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_2.png', 1.232, 0.004, 0.614, 1, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_3.png', 0.972, 0.003, 0.122, 1, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_4.png', 1.218, 0.004, 0.59, 1, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_5.png', 1.284, 0.004, 0.653, 1, 'Rafiqo Rapitasari', 'Benar')
    ,@('AR_1.png', 0.672, 0.002, 0.623, 1, 'Arizli Romadhon', 'Benar')
    ,@('GA_1.png', 1.437, 0.005, 0.61, 1, 'Gege Ardiyansyah', 'Benar')
    ,@('GA_2.png', 0.485, 0.002, 0.308, 0.857, 'Gege Ardiyansyah', 'Benar')
    ,@('GA_3.png', 0.573, 0.002, 0.196, 1, 'Gege Ardiyansyah', 'Benar')
    ,@('FY_1.png', 0.9340000000000001, 0.003, 0.403, 1, 'Fanny Yusuf', 'Benar')
    ,@('FY_2.png', 1.481, 0.005, 0.531, 1, 'Fanny Yusuf', 'Benar')
    ,@('FY_3.png', 1.359, 0.004, 0.539, 0.857, 'Fanny Yusuf', 'Benar')
    ,@('FY_4.png', 1.121, 0.004, 0.428, 0.857, 'Fanny Yusuf', 'Benar')
    ,@('TO_1.png', 0.679, 0.002, 0.441, 1, 'Tiara Oktavian', 'Benar')
    ,@('TO_2.png', 0.901, 0.003, 0.433, 1, 'Tiara Oktavian', 'Benar')
    ,@('TO_3.png', 0.664, 0.002, 0.493, 1, 'Tiara Oktavian', 'Benar')
    ,@('TO_4.png', 2.854, 0.008999999999999999, 0.17, 1, 'Fanny Yusuf', 'Salah')
    ,@('TO_5.png', 2.801, 0.008999999999999999, 0.131, 1, 'Tiara Oktavian', 'Benar')
    ,@('TD_1.png', 2.158, 0.007, 0.035, 0.429, 'Tidak Diketahui', 'Benar')
    ,@('TD_2.png', 2.331, 0.008, 0.061, 0.429, 'Tidak Diketahui', 'Benar')
    ,@('TD_3.png', 0.878, 0.003, 0.211, 0.286, 'Tidak Diketahui', 'Benar')
    ,@('TD_4.png', 0.8, 0.003, 0.109, 0.429, 'Tidak Diketahui', 'Benar')
)

$startRow = 1
$startCol = 1
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, $startCol).Value     = $rec[0]
    $ws.Cells.Item($row, $startCol + 1).Value = $rec[1]
    $ws.Cells.Item($row, $startCol + 2).Value = $rec[2]
    $ws.Cells.Item($row, $startCol + 3).Value = $rec[3]
    $ws.Cells.Item($row, $startCol + 4).Value = $rec[4]
    $ws.Cells.Item($row, $startCol + 5).Value = $rec[5]
    $ws.Cells.Item($row, $startCol + 6).Value = $rec[6]
}
